# Trade #10 closed at 2026-02-17 08:08:34 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: update aggregate stats now that trade #10 has closed.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.96   # Current Capital
$summary.Range("B4").Value = -0.04     # Total P&L $
$summary.Range("B5").Value = -0.08     # Total P&L %
$summary.Range("B6").Value = 10        # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: update MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999   # Capital
$status.Range("D4").Value = 10                  # Trades
$status.Range("E4").Value = -0.04                # P&L $
$status.Range("F4").Value = -0.04                # P&L %
$status.Range("G4").Value = 40                   # Win Rate %

# ---------------------------------------------------------------------
# Append the new closed trade (#10) to both "All Trades" and
# "MarketMaking" sheets as row 11.
# ---------------------------------------------------------------------
$newTradeRow = @(10, "2026-02-17", "08:08:27", "MarketMaking", "UP", 0.33, 0.29, "CLOSED", -12.1212, -0.04, 99.95999999999999, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    # The Date column looks numeric to Excel's auto-detection, so force it
    # to be stored as plain text (matching the other rows) before writing,
    # then drop back to the Normal style so no stray formatting is left
    # behind on the new cell.
    $ws.Cells.Item(11, 2).NumberFormat = "@"
    for ($i = 0; $i -lt $newTradeRow.Length; $i++) {
        $ws.Cells.Item(11, $i + 1).Value = $newTradeRow[$i]
    }
    $ws.Cells.Item(11, 2).Style = "Normal"
}
